# Auto-generated COM-interop script implementing the commit:
# "Append: 2025-09-01 01:28 JST"
#
# Sheet "ランサーズ": insert a fresh scrape result as the new row 2
#   (pushing the previous rows down by one), and rebuild the F-column
#   hyperlinks so they stay aligned with their row.
# Sheet "統計": append one new summary row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ランサーズ")
$ws2 = $wb.Worksheets.Item("統計")

# --- ランサーズ: insert new row 2, then rewrite rows 2..75 in full ---
$ws1.Rows.Item(2).Insert()

# Existing hyperlink refs do not auto-shift with the row insert, so
# drop them all and rebuild from scratch against the final row numbers.
$ws1.Hyperlinks.Delete()

$lancersData = @(
    @("2025-09-01 01:28:36", "【急募】領収書発行システムの構築依頼(簡易的なものでOK)", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383597", 28, $null),
    @("2025-08-31 18:22:49", "急募 【急募】音響系機器のESP32・BP83のファームウェア開発", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383501", 68, "◆開発"),
    @("2025-08-31 18:22:49", "どこにも負けない最速購入", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383527", 13, $null),
    @("2025-08-31 06:23:08", "【急募】小型BLE音声モジュールの試作開発依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383341", 68, "◆開発"),
    @("2025-08-31 06:23:08", "LINE公式アカウントの自動応答・ステップ配信設定(文章提供あり)", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383334", 18, $null),
    @("2025-08-31 01:21:48", "【RPA自動化】通関データ入力業務の効率化依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383211", 153, "◆効率化,自動化"),
    @("2025-08-31 01:21:48", "IB報酬を得るための高性能EA開発依頼", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383199", 68, "◆開発"),
    @("2025-08-31 01:21:48", "マクロの仕様変更をお願いします。", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383131", 10, $null),
    @("2025-08-30 18:22:35", "【急募】AI開発者専門家を探しています!高スキル歓迎", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383100", 368, "🔥AI,Ai ◆開発"),
    @("2025-08-30 18:22:35", "【急募】個人思想を反映したコンテンツ生成システム開発", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383095", 118, "◆開発,システム開発"),
    @("2025-08-30 12:31:43", "サブスク型学習サイトの開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5365024", 93, "◆開発 ◇サイト"),
    @("2025-08-30 12:31:43", "【急募】Google Workspace/AppSheetの設定・開発サポート依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5383065", 60, "◆開発"),
    @("2025-08-30 01:40:08", "【急募】kintone業務アプリ開発のプロフェッショナルを探しています!", "システム開発", "300,000 円 ~", "期限情報なし", "https://www.lancers.jp/work/detail/5359873", 93, "◆開発 ◇アプリ"),
    @("2025-08-30 01:12:26", "AIジャーナリングサービスのWebアプリとLLM開発依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382563", 378, "🔥AI,Ai ◆開発 ◇アプリ"),
    @("2025-08-30 01:12:26", "AIジャーナリングサービスのWebアプリとLLM開発依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382589", 373, "🔥AI,Ai ◆開発 ◇アプリ"),
    @("2025-08-30 01:12:26", "【健康経営支援】体力テストWebアプリ開発のエンジニア募集", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382721", 93, "◆開発 ◇アプリ"),
    @("2025-08-30 01:12:26", "CFDと併用するレンダリング", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382728", 25, $null),
    @("2025-08-30 01:12:26", "EXCELでのマクロ作成", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5379313", 18, $null),
    @("2025-08-30 01:12:26", "【急募】既存WPテーマにWooCommerceでショップページ作成依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382676", 13, $null),
    @("2025-08-30 01:12:26", "【急募】コインパーキング売上データの集計マクロ作成依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5379730", 10, $null),
    @("2025-08-29 18:29:54", "【急募】Next.js × Vercel × Supabaseサイトの公開対応依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5308620", 145, "🔥Next.js ◇サイト"),
    @("2025-08-29 18:29:54", "【業務委託】RPAツールを開発できる方を募集します!", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309123", 128, "◆ツール,開発"),
    @("2025-08-29 18:29:54", "【急募】低レイテンシー設計のリアルタイム監視ツール開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309099", 128, "◆ツール,開発"),
    @("2025-08-29 18:29:54", "WEBアプリの設計及び開発:駐在可能な方", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309115", 100, "◆開発 ◇アプリ"),
    @("2025-08-29 18:29:54", "東海3県の方希望! PHP(ララベル)システムの保守開発 継続長期案件 時々来社希望", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309162", 98, "◆開発 ○PHP"),
    @("2025-08-29 18:29:54", "UWSC自動入力システムの改良・開発者募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5308794", 78, "◆開発"),
    @("2025-08-29 18:29:54", "【急募】ShopifyとPOS連携!会員ランク制度導入開発者募集", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5371807", 68, "◆開発"),
    @("2025-08-29 18:29:54", "東海3県の方希望! EC-CUBEの保守開発 継続長期案件 時々来社打ち合わせ希望", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309166", 60, "◆開発"),
    @("2025-08-29 18:29:54", "初回 【iOS/Androidアプリ再リリース】ランダム通話アプリの修正依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309542", 38, "◇アプリ"),
    @("2025-08-29 18:29:54", "【急募】ワードプレスでのサイト構築依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309147", 33, "◇サイト"),
    @("2025-08-29 18:29:54", "【iOS/Androidアプリ再リリース】ランダム通話アプリの修正依頼", "システム開発", "50,000 円 ~ 60,000 円 / 募集期間 5 日、取引期間 55 日", "期限情報なし", "https://www.lancers.jp/work/detail/5308955", 33, "◇アプリ"),
    @("2025-08-29 18:29:54", "【リプレイス】FileMaker製の業務システムをWeb再構築(受発注・在庫・帳票/CSV/PDF)", "システム開発", "3,000,000 円 ~ 5,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5375452", 40, $null),
    @("2025-08-29 18:29:54", "【継続案件|お気軽にご応募ください】Webシステムのテスター募集", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309229", 25, $null),
    @("2025-08-29 18:29:54", "TradingView自動売買システムの使い方・VPS操作・デプロイ方法を教えてください", "システム開発", "5,000 円 ~ 6,000 円 / 募集期間 3 日、取引期間 2 日", "期限情報なし", "https://www.lancers.jp/work/detail/5308967", 25, $null),
    @("2025-08-29 18:29:54", "初回 【AWS構築】既存のAzureからBubble DBへのデータ移行に必要なAWS中間サーバーの作成", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309519", 18, $null),
    @("2025-08-29 18:29:54", "【フリーランス必見】メモからエッセイを自動生成するボット構築講座", "システム開発", "100,000 円 ~ 200,000 円 / 募集期間 5 日、取引期間 94 日", "期限情報なし", "https://www.lancers.jp/work/detail/5309439", 18, $null),
    @("2025-08-29 18:29:54", "【業務委託】UTAGEを構築できる方を募集します!", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5309131", 18, $null),
    @("2025-08-29 18:24:22", "【急募】フロントエンド Next.js システム開発", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382523", 243, "🔥Next.js ◆開発,システム開発"),
    @("2025-08-29 18:24:22", "【急募】既存Googleスプレッドシートで効率化を図りたいので修正・改善をお手伝いしてくれる方募集", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382629", 70, "◆効率化"),
    @("2025-08-29 18:24:22", "【急募】ZendeskとLINE連携の自動メッセージ設定", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382550", 10, $null),
    @("2025-08-29 12:35:56", "【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382463", 158, "◆自動化,スクレイピング ◇管理"),
    @("2025-08-29 12:35:56", "【急募】大学ウェブサイト修正依頼(WordPress)", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382289", 50, "◇サイト ○WordPress"),
    @("2025-08-29 12:35:56", "【SRE / インフラエンジニア募集(基本リモート・金融系案件)】", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5371075", 25, $null),
    @("2025-08-29 12:35:56", "ESP32を利用した精密機械の製品化サポート", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382344", 10, $null),
    @("2025-08-29 06:24:56", "【急募】心理的安全性を可視化するWebツール開発", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381977", 128, "◆ツール,開発"),
    @("2025-08-29 06:24:56", "【急募】SFAのスクラッチ開発対応してくれる企業様募集!", "システム開発", "5,000,000 円 ~ / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382236", 75, "◆開発"),
    @("2025-08-29 06:24:56", "【ペットのアバター化】LumiGOプロトタイプ開発の依頼", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5382213", 68, "◆開発"),
    @("2025-08-29 01:45:10", "【急募】フロントエンド Webサービス開発/改修", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5370186", 68, "◆開発"),
    @("2025-08-29 01:16:05", "【急募】既存スプレッドシートを新アカウントに移動依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381748", 13, $null),
    @("2025-08-28 12:35:27", "時給2000円|make/n8n/AWSによる自動化開発者を募集(RPA・MCP経験者歓迎)", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381595", 143, "◆開発,自動化"),
    @("2025-08-28 12:35:27", "【急募】愛知県弥富市でAMRとPLCを繋ぐシステム開発", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381608", 125, "◆開発,システム開発"),
    @("2025-08-28 12:35:27", "【急募】ジャストDB(ノーコード)受注システムの編集依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5379679", 33, $null),
    @("2025-08-28 12:35:27", "amazon注文詳細より配送番号と配送業者と配送状況取得してもらいたいです。", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381625", 25, $null),
    @("2025-08-28 12:35:27", "【急募】利用者予定表と勤務表の自動集計をお任せ!", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381634", 13, $null),
    @("2025-08-28 06:27:23", "競馬AIの開発ができる方、もしくはすでに開発済みの方", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381290", 375, "🔥AI,Ai ◆開発"),
    @("2025-08-28 06:27:23", "Google AI studio が生成したウェブアプリの調整", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5371027", 330, "🔥AI,Ai ◇アプリ"),
    @("2025-08-28 06:27:23", "【急募】TEBANASHIの開発補助をお手伝いしてくれる方募集!", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381284", 75, "◆開発"),
    @("2025-08-28 06:27:23", "【急募】MDMサービスの障害対応・追加開発業務", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381264", 75, "◆開発"),
    @("2025-08-28 06:27:23", "【WordPress改善】社内システムの運用支援をお願いします", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381250", 48, "○WordPress"),
    @("2025-08-28 06:27:23", "【急募】会計ソフト導入とVPN構築の専門家を探しています", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381245", 13, $null),
    @("2025-08-28 06:27:23", "【知識ゼロから参加できる♪ 新しい体験プロジェクト】", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381204", 10, $null),
    @("2025-08-28 01:16:20", "【急募】GoogleMAP機能追加・編集の依頼", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5381118", 18, $null),
    @("2025-08-27 18:25:32", "Qt / C++ ベースのWindowsアプリとデバイスファームウェアの修正", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380896", 45, "◇アプリ"),
    @("2025-08-27 12:36:04", "【急募】outsystemsでホテル予約アプリの開発依頼", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380830", 85, "◆開発 ◇アプリ"),
    @("2025-08-27 12:36:04", "wordpressレンダリングを妨げるリソースの除外", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5016989", 33, "○WordPress"),
    @("2025-08-27 06:56:02", "Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5273634", 298, "🔥Python ◆開発,スクレイピング"),
    @("2025-08-27 06:56:02", "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5217096", 243, "🔥API ◆ツール"),
    @("2025-08-27 06:56:02", "【急募】youtube動画配信システム開発のフリーランス募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380343", 118, "◆開発,システム開発"),
    @("2025-08-27 06:56:02", "フロントエンドエンジニア マッチングサイトの開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380337", 93, "◆開発 ◇サイト"),
    @("2025-08-27 06:56:02", "【急募】PHPサイトのDB情報をルート外に移動", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380683", 50, "◇サイト ○PHP"),
    @("2025-08-27 06:56:02", "Bubbleアプリ機能改善依頼書", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380747", 33, "◇アプリ"),
    @("2025-08-27 06:56:02", "〖リモート可〗Delphiエンジニア募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5341051", 25, $null),
    @("2025-08-27 06:56:02", "【急募】シティヘブン自動予約システム構築の依頼", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380357", 25, $null),
    @("2025-08-27 06:56:02", "初回 YAMAHAルーターRTX830で自社サーバー1台のポート開放設定をお願いします。", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5380420", 13, $null)
)

for ($i = 0; $i -lt $lancersData.Length; $i++) {
    $row = $lancersData[$i]
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Hyperlinks.Add($ws1.Cells.Item($r, 6), $row[5]) | Out-Null
    $ws1.Cells.Item($r, 6).Style = "Hyperlink"
    $ws1.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne $null) {
        $ws1.Cells.Item($r, 8).Value = $row[7]
    }
}

# --- 統計: append the new summary row (row 39) ---
$statRow = @("2025-09-01T01:28:36.042462", 9, "全案件リスト", 55.6, 4, 2, 9)
$statR = $ws2.UsedRange.Rows.Count + 1
$ws2.Cells.Item($statR, 1).Value = $statRow[0]
$ws2.Cells.Item($statR, 2).Value = $statRow[1]
$ws2.Cells.Item($statR, 3).Value = $statRow[2]
$ws2.Cells.Item($statR, 4).Value = $statRow[3]
$ws2.Cells.Item($statR, 5).Value = $statRow[4]
$ws2.Cells.Item($statR, 6).Value = $statRow[5]
$ws2.Cells.Item($statR, 7).Value = $statRow[6]

